$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7866.433
$ws.Range("I132").Value = 7999.75
$ws.Range("K132").Value = 23999.25
$ws.Range("M132").Value = -21469.25
$ws.Range("H137").Value = 26325102
$ws.Range("I137").Value = 38463264
$ws.Range("K137").Value = 115389792
$ws.Range("M137").Value = -115387242
$ws.Range("H138").Value = 4333
$ws.Range("I138").Value = 4651.091
$ws.Range("K138").Value = 13953.273
$ws.Range("M138").Value = -8813.273000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1428992.4
$ws.Range("I32").Value = 1544929.2
$ws.Range("J32").Value = 37750
$ws.Range("K32").Value = 1544929.2
$ws.Range("L32").Value = 37750
$ws.Range("M32").Value = -1544642.2
$ws.Range("N32").Value = -38324
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H74").Value = 3282961.8
$ws.Range("I74").Value = 4284103
$ws.Range("J74").Value = 29253.5
$ws.Range("K74").Value = 4284103
$ws.Range("L74").Value = 29253.5
$ws.Range("M74").Value = -4283229
$ws.Range("N74").Value = -31001.5
$ws.Range("H77").Value = 3282961.8
$ws.Range("I77").Value = 4284103
$ws.Range("J77").Value = 29253.5
$ws.Range("K77").Value = 21420515
$ws.Range("L77").Value = 146267.5
$ws.Range("M77").Value = -21416147
$ws.Range("N77").Value = -155003.5
$ws.Range("H102").Value = 1543.4445
$ws.Range("I102").Value = 1486.375
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1486.375
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 135.625
$ws.Range("N102").Value = -5244
$ws.Range("H110").Value = 1885.5883
$ws.Range("I110").Value = 1722.2858
$ws.Range("J110").Value = 1999.9
$ws.Range("K110").Value = 1722.2858
$ws.Range("L110").Value = 1999.9
$ws.Range("M110").Value = 322.7141999999999
$ws.Range("N110").Value = -6089.9
$ws.Range("H132").Value = 4575.7837
$ws.Range("I132").Value = 3562.3572
$ws.Range("K132").Value = 10687.0716
$ws.Range("M132").Value = -8157.071599999999
$ws.Range("H134").Value = 52500.5
$ws.Range("J134").Value = 52500.5
$ws.Range("L134").Value = 52500.5
$ws.Range("N134").Value = -62640.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1718.421
$ws.Range("I107").Value = 1405.5555
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1405.5555
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 514.4445000000001
$ws.Range("N107").Value = -5840
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").ClearContents()
$ws.Range("N131").Value = 0
$ws.Range("H135").Value = 52033.715
$ws.Range("J135").Value = 52033.715
$ws.Range("L135").Value = 52033.715
$ws.Range("N135").Value = -62173.715

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 12000
$ws.Range("I69").Value = 12000
$ws.Range("K69").Value = 12000
$ws.Range("M69").Value = -11251
$ws.Range("H72").Value = 12000
$ws.Range("I72").Value = 12000
$ws.Range("K72").Value = 36000
$ws.Range("M72").Value = -32256
$ws.Range("H134").Value = 3738.9546
$ws.Range("I134").Value = 3738.9546
$ws.Range("K134").Value = 11216.8638
$ws.Range("M134").Value = -8681.863799999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").ClearContents()
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = 0
$ws.Range("H34").Value = 4399.5557
$ws.Range("J34").Value = 7944.4443
$ws.Range("L34").Value = 23833.3329
$ws.Range("N34").Value = -24001.3329
$ws.Range("H80").Value = 10562.625
$ws.Range("J80").Value = 11643
$ws.Range("L80").Value = 34929
$ws.Range("N80").Value = -36801
$ws.Range("H83").Value = 10562.625
$ws.Range("J83").Value = 11643
$ws.Range("L83").Value = 104787
$ws.Range("N83").Value = -114147
$ws.Range("H132").Value = 2271.5715
$ws.Range("J132").Value = 4212
$ws.Range("L132").Value = 37908
$ws.Range("N132").Value = -42968
$ws.Range("H140").Value = 2395.5
$ws.Range("I140").Value = 2395.5
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 7186.5
$ws.Range("L140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -2006.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1341.6666
$ws.Range("I107").Value = 1900
$ws.Range("K107").Value = 1900
$ws.Range("M107").Value = 20
$ws.Range("H126").Value = 3357.2
$ws.Range("I126").Value = 3143
$ws.Range("K126").Value = 9429
$ws.Range("M126").Value = -6959
$ws.Range("H132").Value = 10050.16
$ws.Range("I132").Value = 6952.75
$ws.Range("K132").Value = 20858.25
$ws.Range("M132").Value = -18328.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4205.375
$ws.Range("I7").Value = 4446.6665
$ws.Range("J7").Value = 3481.5
$ws.Range("K7").Value = 4446.6665
$ws.Range("L7").Value = 3481.5
$ws.Range("M7").Value = -4334.6665
$ws.Range("N7").Value = -3705.5
$ws.Range("H43").Value = 180166.67
$ws.Range("I43").Value = 18000
$ws.Range("J43").Value = 200437.5
$ws.Range("K43").Value = 18000
$ws.Range("L43").Value = 200437.5
$ws.Range("M43").Value = -17807
$ws.Range("N43").Value = -200823.5
$ws.Range("H54").Value = 36646
$ws.Range("J54").Value = 36646
$ws.Range("L54").Value = 36646
$ws.Range("N54").Value = -37934
$ws.Range("H55").Value = 1378.5883
$ws.Range("I55").Value = 1307.6154
$ws.Range("J55").Value = 1422.5238
$ws.Range("K55").Value = 1307.6154
$ws.Range("L55").Value = 1422.5238
$ws.Range("M55").Value = -1134.6154
$ws.Range("N55").Value = -1768.5238
$ws.Range("H126").Value = 4205.375
$ws.Range("I126").Value = 4446.6665
$ws.Range("J126").Value = 3481.5
$ws.Range("K126").Value = 13339.9995
$ws.Range("L126").Value = 10444.5
$ws.Range("M126").Value = -10869.9995
$ws.Range("N126").Value = -15384.5
$ws.Range("H132").Value = 1014464.9
$ws.Range("I132").Value = 2225237.5
$ws.Range("K132").Value = 6675712.5
$ws.Range("M132").Value = -6673182.5
$ws.Range("H136").Value = 10872264
$ws.Range("I136").Value = 7355444
$ws.Range("K136").Value = 22066332
$ws.Range("M136").Value = -22063782

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 75979.47
$ws.Range("J81").Value = 124921.22
$ws.Range("L81").Value = 249842.44
$ws.Range("N81").Value = -251964.44
$ws.Range("H84").Value = 75979.47
$ws.Range("J84").Value = 124921.22
$ws.Range("L84").Value = 1249212.2
$ws.Range("N84").Value = -1259820.2
$ws.Range("H132").Value = 4903798
$ws.Range("I132").Value = 6412005
$ws.Range("J132").Value = 2124
$ws.Range("K132").Value = 19236015
$ws.Range("L132").Value = 6372
$ws.Range("M132").Value = -19233485
$ws.Range("N132").Value = -11432
$ws.Range("H133").Value = 51000
$ws.Range("J133").Value = 51000
$ws.Range("L133").Value = 51000
$ws.Range("N133").Value = -61120
$ws.Range("H136").Value = 10599344
$ws.Range("I136").Value = 6212364.5
$ws.Range("J136").Value = 20835628
$ws.Range("K136").Value = 18637093.5
$ws.Range("L136").Value = 62506884
$ws.Range("M136").Value = -18634543.5
$ws.Range("N136").Value = -62511984
